$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "June 28, 2023"
$ws.Range("E2").Value = 21064
$ws.Range("F2").Value = 15682
$ws.Range("G2").Value = 3192
$ws.Range("H2").Value = 690
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 1350
$ws.Range("B3").Value = "June 28, 2023"
$ws.Range("E3").Value = 317
$ws.Range("F3").Value = 225
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 36
$ws.Range("B4").Value = "June 28, 2023"
$ws.Range("B5").Value = "June 28, 2023"
$ws.Range("E5").Value = 90385
$ws.Range("F5").Value = 67964
$ws.Range("G5").Value = 10939
$ws.Range("H5").Value = 3986
$ws.Range("I5").Value = 1517
$ws.Range("J5").Value = 5979
$ws.Range("B6").Value = "June 28, 2023"
$ws.Range("E6").Value = 996736
$ws.Range("F6").Value = 770514
$ws.Range("G6").Value = 153002
$ws.Range("H6").Value = 29950
$ws.Range("I6").Value = 7376
$ws.Range("J6").Value = 35894
$ws.Range("B7").Value = "June 28, 2023"
$ws.Range("E7").Value = 190393
$ws.Range("F7").Value = 128914
$ws.Range("G7").Value = 18822
$ws.Range("H7").Value = 8628
$ws.Range("I7").Value = 6101
$ws.Range("J7").Value = 27928
$ws.Range("B8").Value = "June 28, 2023"
$ws.Range("E8").Value = 59644
$ws.Range("F8").Value = 33918
$ws.Range("G8").Value = 4664
$ws.Range("H8").Value = 2494
$ws.Range("I8").Value = 2816
$ws.Range("J8").Value = 15752
$ws.Range("B9").Value = "June 28, 2023"
$ws.Range("E9").Value = 774603
$ws.Range("F9").Value = 425981
$ws.Range("G9").Value = 71222
$ws.Range("H9").Value = 42219
$ws.Range("I9").Value = 31260
$ws.Range("J9").Value = 203921
$ws.Range("B10").Value = "June 28, 2023"
$ws.Range("E10").Value = 12614
$ws.Range("F10").Value = 5556
$ws.Range("G10").Value = 914
$ws.Range("H10").Value = 863
$ws.Range("I10").Value = 852
$ws.Range("J10").Value = 4429
$ws.Range("B11").Value = "June 28, 2023"
$ws.Range("E11").Value = 5395
$ws.Range("F11").Value = 3241
$ws.Range("G11").Value = 792
$ws.Range("J11").Value = 682
$ws.Range("B12").Value = "June 28, 2023"
$ws.Range("E12").Value = 636602
$ws.Range("F12").Value = 482279
$ws.Range("G12").Value = 106609
$ws.Range("H12").Value = 29928
$ws.Range("I12").Value = 4861
$ws.Range("J12").Value = 12925
$ws.Range("B13").Value = "June 28, 2023"
$ws.Range("E13").Value = 143479
$ws.Range("F13").Value = 92943
$ws.Range("H13").Value = 7831
$ws.Range("J13").Value = 18949
$ws.Range("B14").Value = "June 28, 2023"
$ws.Range("E14").Value = 76809
$ws.Range("F14").Value = 54799
$ws.Range("G14").Value = 9141
$ws.Range("H14").Value = 3698
$ws.Range("I14").Value = 1744
$ws.Range("J14").Value = 7427
$ws.Range("B15").Value = "June 28, 2023"
$ws.Range("E15").Value = 55377
$ws.Range("F15").Value = 34525
$ws.Range("G15").Value = 7764
$ws.Range("H15").Value = 3075
$ws.Range("I15").Value = 1692
$ws.Range("J15").Value = 8321
$ws.Range("B16").Value = "June 28, 2023"
$ws.Range("E16").Value = 29638
$ws.Range("F16").Value = 19405
$ws.Range("G16").Value = 4301
$ws.Range("H16").Value = 1742
$ws.Range("I16").Value = 945
$ws.Range("J16").Value = 3245
$ws.Range("B17").Value = "June 28, 2023"
$ws.Range("E17").Value = 253
$ws.Range("F17").Value = 194
$ws.Range("B18").Value = "June 28, 2023"
$ws.Range("E18").Value = 189413
$ws.Range("F18").Value = 146028
$ws.Range("G18").Value = 28299
$ws.Range("H18").Value = 7675
$ws.Range("I18").Value = 1766
$ws.Range("J18").Value = 5645
$ws.Range("B19").Value = "June 28, 2023"
$ws.Range("E19").Value = 277555
$ws.Range("F19").Value = 106060
$ws.Range("G19").Value = 20838
$ws.Range("H19").Value = 16279
$ws.Range("I19").Value = 19003
$ws.Range("J19").Value = 115375
$ws.Range("B20").Value = "June 28, 2023"
$ws.Range("B21").Value = "June 28, 2023"
$ws.Range("B22").Value = "June 28, 2023"
$ws.Range("E22").Value = 379702
$ws.Range("F22").Value = 247736
$ws.Range("G22").Value = 36011
$ws.Range("H22").Value = 20461
$ws.Range("I22").Value = 12282
$ws.Range("J22").Value = 63212
$ws.Range("B23").Value = "June 28, 2023"
$ws.Range("E23").Value = 36372
$ws.Range("F23").Value = 24035
$ws.Range("G23").Value = 4776
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 974
$ws.Range("J23").Value = 4587
$ws.Range("B24").Value = "June 28, 2023"
$ws.Range("E24").Value = 153847
$ws.Range("F24").Value = 80686
$ws.Range("G24").Value = 15374
$ws.Range("H24").Value = 9170
$ws.Range("I24").Value = 7775
$ws.Range("J24").Value = 40842
$ws.Range("B25").Value = "June 28, 2023"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.8"
$ws.Range("F25").Value = 684
$ws.Range("G25").Value = 89
$ws.Range("H25").Value = 79
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 241
$ws.Range("B26").Value = "June 28, 2023"
$ws.Range("E26").Value = 2953
$ws.Range("J26").Value = 758
$ws.Range("B27").Value = "June 28, 2023"
$ws.Range("E27").Value = 461
$ws.Range("F27").Value = 287
$ws.Range("G27").Value = 74
$ws.Range("B28").Value = "June 28, 2023"
$ws.Range("B29").Value = "June 28, 2023"
$ws.Range("E29").Value = 1540
$ws.Range("F29").Value = 533
$ws.Range("J29").Value = 721
$ws.Range("B30").Value = "June 28, 2023"
$ws.Range("B31").Value = "June 28, 2023"
$ws.Range("B32").Value = "June 28, 2023"
$ws.Range("B33").Value = "June 28, 2023"
$ws.Range("E33").Value = 20534
$ws.Range("F33").Value = 12695
$ws.Range("G33").Value = 2939
$ws.Range("H33").Value = 1826
$ws.Range("I33").Value = 614
$ws.Range("J33").Value = 2460
$ws.Range("B34").Value = "June 28, 2023"
$ws.Range("B35").Value = "June 28, 2023"
$ws.Range("B36").Value = "June 28, 2023"
$ws.Range("B37").Value = "June 28, 2023"
$ws.Range("E37").Value = 884
$ws.Range("F37").Value = 646
$ws.Range("G37").Value = 97
$ws.Range("H37").Value = 26
$ws.Range("I37").Value = 53
$ws.Range("J37").Value = 62
$ws.Range("B38").Value = "June 28, 2023"
